# Auto-generated PowerShell COM-interop script
# Applies numeric 'want to go' count updates and one address text fix
# across sheets 展览, 演出, 本地生活, 全部类型 per commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 75
$ws.Range("D5").Value = "共和新路3201号 静安国际科创社区云芯科创中心"
$ws.Range("F5").Value = 785
$ws.Range("F6").Value = 121
$ws.Range("F7").Value = 516
$ws.Range("F8").Value = 958
$ws.Range("F9").Value = 1643
$ws.Range("F10").Value = 1297
$ws.Range("F11").Value = 1583
$ws.Range("F13").Value = 1581
$ws.Range("F14").Value = 352
$ws.Range("F15").Value = 1714
$ws.Range("F17").Value = 1142
$ws.Range("F18").Value = 385
$ws.Range("F19").Value = 58
$ws.Range("F20").Value = 117
$ws.Range("F21").Value = 1914
$ws.Range("F22").Value = 264
$ws.Range("F24").Value = 1016
$ws.Range("F26").Value = 1291
$ws.Range("F27").Value = 1087
$ws.Range("F28").Value = 89
$ws.Range("F30").Value = 1215
$ws.Range("F32").Value = 1194
$ws.Range("F33").Value = 1150
$ws.Range("F34").Value = 295
$ws.Range("F35").Value = 89
$ws.Range("F36").Value = 900
$ws.Range("F38").Value = 1717
$ws.Range("F40").Value = 125
$ws.Range("F41").Value = 2089
$ws.Range("F42").Value = 102
$ws.Range("F43").Value = 844
$ws.Range("F44").Value = 105
$ws.Range("F47").Value = 124

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1509
$ws.Range("F8").Value = 2616
$ws.Range("F12").Value = 268
$ws.Range("F21").Value = 328
$ws.Range("F26").Value = 38
$ws.Range("F27").Value = 38
$ws.Range("F33").Value = 57
$ws.Range("F35").Value = 32
$ws.Range("F42").Value = 69

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2937
$ws.Range("F6").Value = 4714
$ws.Range("F9").Value = 604
$ws.Range("F10").Value = 795
$ws.Range("F11").Value = 487
$ws.Range("F12").Value = 437
$ws.Range("F13").Value = 1173
$ws.Range("F14").Value = 326
$ws.Range("F15").Value = 802

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 795
$ws.Range("F6").Value = 75
$ws.Range("F7").Value = 437
$ws.Range("F8").Value = 437
$ws.Range("F9").Value = 1173
$ws.Range("F10").Value = 516
$ws.Range("F11").Value = 958
$ws.Range("F12").Value = 1644
$ws.Range("F13").Value = 1297
$ws.Range("F14").Value = 1583
$ws.Range("F16").Value = 1581
$ws.Range("F17").Value = 268
$ws.Range("F19").Value = 1714
$ws.Range("F20").Value = 1142
$ws.Range("F22").Value = 802
$ws.Range("F23").Value = 802
$ws.Range("F24").Value = 1915
$ws.Range("F25").Value = 264
$ws.Range("F27").Value = 1016
$ws.Range("F29").Value = 1291
$ws.Range("F30").Value = 328
$ws.Range("F31").Value = 1087
$ws.Range("F32").Value = 89
$ws.Range("F33").Value = 1215
$ws.Range("F35").Value = 1194
$ws.Range("F37").Value = 38
$ws.Range("F38").Value = 1151
$ws.Range("F39").Value = 295
$ws.Range("F40").Value = 900
$ws.Range("F43").Value = 1717
$ws.Range("F45").Value = 125
$ws.Range("F46").Value = 2089
$ws.Range("F47").Value = 102
$ws.Range("F48").Value = 844
$ws.Range("F49").Value = 105
$ws.Range("F51").Value = 124
